$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C/A Lag row's C value: 2.51 -> 2.546 (kept as text, matching original type)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2.546"
$ws.Range("C2").Style = "Normal"

# Remove the "Constant" and "r2_adj" rows entirely
$ws.Range("A4:C5").Delete()
